$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 123
$ws1.Range("F3").Value = 229
$ws1.Range("F5").Value = 6643
$ws1.Range("F9").Value = 6083
$ws1.Range("F14").Value = 90
$ws1.Range("F15").Value = 387
$ws1.Range("F21").Value = 4393
$ws1.Range("F23").Value = 20
$ws1.Range("F25").Value = 35

# Sheet "全部类型" (Worksheets index 4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 123
$ws4.Range("F3").Value = 229
$ws4.Range("F5").Value = 6643
$ws4.Range("F9").Value = 6083
$ws4.Range("F14").Value = 90
$ws4.Range("F15").Value = 387
$ws4.Range("F21").Value = 4393
$ws4.Range("F24").Value = 20
$ws4.Range("F26").Value = 35
